$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (shifts existing rows 5-18 down to 6-19)
$ws.Rows.Item(5).Insert()

# Fill in the new row 5 with the Aviation Gasoline mapping entry.
# Set C5 before A5 so the new shared strings are appended in the same
# order as the target workbook (subfuel code string before the label).
$ws.Range("B5").Value = "07_petroleum_products"
$ws.Range("C5").Value = "07_02_aviation_gasoline"
$ws.Range("A5").Value = "Aviation Gasoline"
